# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect newly scraped counts (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Row -> new value updates for the "展览" sheet (column F)
$wsExhibit.Range("F3").Value = 246
$wsExhibit.Range("F5").Value = 2868
$wsExhibit.Range("F12").Value = 82
$wsExhibit.Range("F16").Value = 4704
$wsExhibit.Range("F18").Value = 5117
$wsExhibit.Range("F19").Value = 1674
$wsExhibit.Range("F20").Value = 2871
$wsExhibit.Range("F21").Value = 3268
$wsExhibit.Range("F29").Value = 1872
$wsExhibit.Range("F32").Value = 710

# Matching row -> new value updates for the "全部类型" sheet (column F)
$wsAll.Range("F8").Value = 246
$wsAll.Range("F11").Value = 2868
$wsAll.Range("F19").Value = 82
$wsAll.Range("F26").Value = 4704
$wsAll.Range("F28").Value = 5117
$wsAll.Range("F29").Value = 1674
$wsAll.Range("F30").Value = 2871
$wsAll.Range("F31").Value = 3268
$wsAll.Range("F43").Value = 1872
$wsAll.Range("F46").Value = 710
